# Auto-generated edit script: updates FFXIV leve-profit market price data
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per refreshed
# market-board snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1771.8334
$ws.Range("I11").Value = 1771.8334
$ws.Range("K11").Value = 1771.8334
$ws.Range("M11").Value = -1631.8334

$ws.Range("H15").Value = 1416.7407
$ws.Range("I15").Value = 1416.7407
$ws.Range("K15").Value = 4250.2221
$ws.Range("M15").Value = -4081.2221

$ws.Range("H18").Value = 2198.1428
$ws.Range("I18").Value = 2197.8333
$ws.Range("K18").Value = 2197.8333
$ws.Range("M18").Value = -1913.8333

$ws.Range("H40").Value = 4998.6665
$ws.Range("I40").Value = 4998.6665
$ws.Range("K40").Value = 4998.6665
$ws.Range("M40").Value = -4823.6665

$ws.Range("H51").Value = 29541.455
$ws.Range("J51").Value = 9992
$ws.Range("L51").Value = 9992
$ws.Range("N51").Value = -10960

$ws.Range("H64").Value = 3749.5
$ws.Range("I64").Value = 3749.5
$ws.Range("K64").Value = 3749.5
$ws.Range("M64").Value = -3501.5

$ws.Range("H67").Value = 3749.5
$ws.Range("I67").Value = 3749.5
$ws.Range("K67").Value = 3749.5
$ws.Range("M67").Value = -2891.5

$ws.Range("H74").Value = 6908.12
$ws.Range("I74").Value = 6357
$ws.Range("J74").Value = 7167.4707
$ws.Range("K74").Value = 6357
$ws.Range("L74").Value = 7167.4707
$ws.Range("M74").Value = -5421
$ws.Range("N74").Value = -9039.4707

$ws.Range("H77").Value = 6908.12
$ws.Range("I77").Value = 6357
$ws.Range("J77").Value = 7167.4707
$ws.Range("K77").Value = 31785
$ws.Range("L77").Value = 35837.3535
$ws.Range("M77").Value = -27105
$ws.Range("N77").Value = -45197.3535

$ws.Range("H88").Value = 1486.4615
$ws.Range("I88").Value = 1616.6666
$ws.Range("J88").Value = 1447.4
$ws.Range("K88").Value = 1616.6666
$ws.Range("L88").Value = 1447.4
$ws.Range("M88").Value = -1210.6666
$ws.Range("N88").Value = -2259.4

$ws.Range("H91").Value = 1486.4615
$ws.Range("I91").Value = 1616.6666
$ws.Range("J91").Value = 1447.4
$ws.Range("K91").Value = 1616.6666
$ws.Range("L91").Value = 1447.4
$ws.Range("M91").Value = -212.6666
$ws.Range("N91").Value = -4255.4

$ws.Range("H100").Value = 1302.1904
$ws.Range("I100").Value = 972.3125
$ws.Range("K100").Value = 972.3125
$ws.Range("M100").Value = -431.3125

$ws.Range("H107").Value = 1745.238
$ws.Range("I107").Value = 1811.9412
$ws.Range("J107").Value = 1461.75
$ws.Range("K107").Value = 1811.9412
$ws.Range("L107").Value = 1461.75
$ws.Range("M107").Value = 108.0588
$ws.Range("N107").Value = -5301.75

$ws.Range("H113").Value = 6234.375
$ws.Range("I113").Value = 5281.4165
$ws.Range("J113").Value = 7187.3335
$ws.Range("K113").Value = 5281.4165
$ws.Range("L113").Value = 7187.3335
$ws.Range("M113").Value = -2027.4165
$ws.Range("N113").Value = -13695.3335

$ws.Range("H115").Value = 471.33334
$ws.Range("I115").Value = 471.33334
$ws.Range("K115").Value = 1414.00002
$ws.Range("M115").Value = 152.9999800000001

$ws.Range("H116").Value = 4527.6665
$ws.Range("J116").Value = 3246.5
$ws.Range("L116").Value = 3246.5
$ws.Range("N116").Value = -10130.5

$ws.Range("H137").Value = 2521.139
$ws.Range("I137").Value = 1490.1904
$ws.Range("K137").Value = 4470.5712
$ws.Range("M137").Value = -1920.5712

$ws.Range("H138").Value = 4161.986
$ws.Range("J138").Value = 3990.3845
$ws.Range("L138").Value = 11971.1535
$ws.Range("N138").Value = -22251.1535

$ws.Range("H141").Value = 8610.615
$ws.Range("I141").Value = 4471
$ws.Range("K141").Value = 13413
$ws.Range("M141").Value = -8233

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1450.9412
$ws.Range("I2").Value = 1361.7858
$ws.Range("J2").Value = 1867
$ws.Range("K2").Value = 1361.7858
$ws.Range("L2").Value = 1867
$ws.Range("M2").Value = -1248.7858
$ws.Range("N2").Value = -2093

$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").ClearContents()

$ws.Range("H44").Value = 5000
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -4512
$ws.Range("N44").ClearContents()

$ws.Range("H63").Value = 1638.75
$ws.Range("I63").Value = 1444.2858
$ws.Range("K63").Value = 1444.2858
$ws.Range("M63").Value = -758.2858000000001

$ws.Range("H66").Value = 1638.75
$ws.Range("I66").Value = 1444.2858
$ws.Range("K66").Value = 7221.429
$ws.Range("M66").Value = -3789.429

$ws.Range("H88").Value = 15363.5
$ws.Range("J88").Value = 17129.715
$ws.Range("L88").Value = 17129.715
$ws.Range("N88").Value = -17941.715

$ws.Range("H91").Value = 15363.5
$ws.Range("J91").Value = 17129.715
$ws.Range("L91").Value = 17129.715
$ws.Range("N91").Value = -19937.715

$ws.Range("H94").Value = 27665
$ws.Range("J94").Value = 27665
$ws.Range("L94").Value = 27665
$ws.Range("N94").Value = -29467

$ws.Range("H103").Value = 40786.668
$ws.Range("J103").Value = 40786.668
$ws.Range("L103").Value = 40786.668
$ws.Range("N103").Value = -43130.668

$ws.Range("H110").Value = 1057.7778
$ws.Range("I110").Value = 882.5
$ws.Range("J110").Value = 2460
$ws.Range("K110").Value = 882.5
$ws.Range("L110").Value = 2460
$ws.Range("M110").Value = 1162.5
$ws.Range("N110").Value = -6550

$ws.Range("H116").Value = 1450.9412
$ws.Range("I116").Value = 1361.7858
$ws.Range("J116").Value = 1867
$ws.Range("K116").Value = 1361.7858
$ws.Range("L116").Value = 1867
$ws.Range("M116").Value = 932.2141999999999
$ws.Range("N116").Value = -6455

$ws.Range("H132").Value = 4339.515
$ws.Range("I132").Value = 3852.0386
$ws.Range("K132").Value = 11556.1158
$ws.Range("M132").Value = -9026.1158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2385.7273
$ws.Range("I3").Value = 520.7143
$ws.Range("J3").Value = 5649.5
$ws.Range("K3").Value = 520.7143
$ws.Range("L3").Value = 5649.5
$ws.Range("M3").Value = -381.7143
$ws.Range("N3").Value = -5927.5

$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("N4").ClearContents()

$ws.Range("H64").Value = 325.86667
$ws.Range("I64").Value = 268.25
$ws.Range("J64").Value = 346.81818
$ws.Range("K64").Value = 268.25
$ws.Range("L64").Value = 346.81818
$ws.Range("M64").Value = -43.25
$ws.Range("N64").Value = -796.81818

$ws.Range("H67").Value = 325.86667
$ws.Range("I67").Value = 268.25
$ws.Range("J67").Value = 346.81818
$ws.Range("K67").Value = 268.25
$ws.Range("L67").Value = 346.81818
$ws.Range("M67").Value = 511.75
$ws.Range("N67").Value = -1906.81818

$ws.Range("H74").Value = 30460
$ws.Range("J74").Value = 30460
$ws.Range("L74").Value = 30460
$ws.Range("N74").Value = -32332

$ws.Range("H77").Value = 30460
$ws.Range("J77").Value = 30460
$ws.Range("L77").Value = 91380
$ws.Range("N77").Value = -100740

$ws.Range("H86").Value = 9426.9375
$ws.Range("I86").Value = 2827.9092
$ws.Range("K86").Value = 2827.9092
$ws.Range("M86").Value = -1704.9092

$ws.Range("H89").Value = 9426.9375
$ws.Range("I89").Value = 2827.9092
$ws.Range("K89").Value = 14139.546
$ws.Range("M89").Value = -8523.546

$ws.Range("H105").Value = 2270.625
$ws.Range("I105").Value = 2336.85
$ws.Range("K105").Value = 2336.85
$ws.Range("M105").Value = -589.8499999999999

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H117").Value = 118999.664
$ws.Range("J117").Value = 118999.664
$ws.Range("L117").Value = 118999.664
$ws.Range("N117").Value = -128177.664

$ws.Range("H134").Value = 28129.572
$ws.Range("I134").Value = 42637.332
$ws.Range("J134").Value = 17248.75
$ws.Range("K134").Value = 127911.996
$ws.Range("L134").Value = 51746.25
$ws.Range("M134").Value = -125376.996
$ws.Range("N134").Value = -56816.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2385.7273
$ws.Range("I10").Value = 520.7143
$ws.Range("J10").Value = 5649.5
$ws.Range("K10").Value = 520.7143
$ws.Range("L10").Value = 5649.5
$ws.Range("M10").Value = -381.7143
$ws.Range("N10").Value = -5927.5

$ws.Range("H12").Value = 3000
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = -1830
$ws.Range("N12").Value = -4340

$ws.Range("H16").Value = 2998.5
$ws.Range("I16").Value = 1250.5
$ws.Range("J16").Value = 6494.5
$ws.Range("K16").Value = 1250.5
$ws.Range("L16").Value = 6494.5
$ws.Range("M16").Value = -963.5
$ws.Range("N16").Value = -7068.5

$ws.Range("H31").Value = 2231.566
$ws.Range("I31").Value = 1603.7084
$ws.Range("J31").Value = 2751.1724
$ws.Range("K31").Value = 1603.7084
$ws.Range("L31").Value = 2751.1724
$ws.Range("M31").Value = -1308.7084
$ws.Range("N31").Value = -3341.1724

$ws.Range("H34").Value = 2231.566
$ws.Range("I34").Value = 1603.7084
$ws.Range("J34").Value = 2751.1724
$ws.Range("K34").Value = 1603.7084
$ws.Range("L34").Value = 2751.1724
$ws.Range("M34").Value = -1401.7084
$ws.Range("N34").Value = -3155.1724

$ws.Range("H58").Value = 3800
$ws.Range("I58").Value = 1834
$ws.Range("K58").Value = 1834
$ws.Range("M58").Value = -1631

$ws.Range("H94").Value = 2864.5
$ws.Range("I94").Value = 3144.5
$ws.Range("J94").Value = 2771.1667
$ws.Range("K94").Value = 3144.5
$ws.Range("L94").Value = 2771.1667
$ws.Range("M94").Value = -2693.5
$ws.Range("N94").Value = -3673.1667

$ws.Range("H105").Value = 1493.0435
$ws.Range("I105").Value = 1226.2941
$ws.Range("K105").Value = 1226.2941
$ws.Range("M105").Value = 520.7058999999999

$ws.Range("H113").Value = 2998.5
$ws.Range("I113").Value = 1250.5
$ws.Range("J113").Value = 6494.5
$ws.Range("K113").Value = 1250.5
$ws.Range("L113").Value = 6494.5
$ws.Range("M113").Value = 919.5
$ws.Range("N113").Value = -10834.5

$ws.Range("H122").Value = 2804.0908
$ws.Range("I122").Value = 2730.875
$ws.Range("J122").Value = 2999.3333
$ws.Range("K122").Value = 8192.625
$ws.Range("L122").Value = 8997.999899999999
$ws.Range("M122").Value = -5742.625
$ws.Range("N122").Value = -13897.9999

$ws.Range("H136").Value = 3800
$ws.Range("I136").Value = 1834
$ws.Range("K136").Value = 5502
$ws.Range("M136").Value = -2952

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1104.8695
$ws.Range("I107").Value = 495.3
$ws.Range("J107").Value = 1573.7693
$ws.Range("K107").Value = 1485.9
$ws.Range("L107").Value = 4721.3079
$ws.Range("M107").Value = 434.0999999999999
$ws.Range("N107").Value = -8561.3079

$ws.Range("H124").Value = 18283.5
$ws.Range("I124").Value = 2999
$ws.Range("K124").Value = 8997
$ws.Range("M124").Value = -4087

$ws.Range("H127").Value = 433
$ws.Range("J127").Value = 433
$ws.Range("L127").Value = 1299
$ws.Range("N127").Value = -11219

$ws.Range("H131").Value = 3826.1155
$ws.Range("J131").Value = 3810.87
$ws.Range("L131").Value = 11432.61
$ws.Range("N131").Value = -21512.61

$ws.Range("H132").Value = 7412.467
$ws.Range("I132").Value = 6191
$ws.Range("K132").Value = 55719
$ws.Range("M132").Value = -53189

$ws.Range("H141").Value = 34288.91
$ws.Range("I141").Value = 4970.5
$ws.Range("J141").Value = 45283.312
$ws.Range("K141").Value = 14911.5
$ws.Range("L141").Value = 135849.936
$ws.Range("M141").Value = -9731.5
$ws.Range("N141").Value = -146209.936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2299
$ws.Range("I80").Value = 2165
$ws.Range("K80").Value = 2165
$ws.Range("M80").Value = -1167

$ws.Range("H83").Value = 2299
$ws.Range("I83").Value = 2165
$ws.Range("K83").Value = 10825
$ws.Range("M83").Value = -5833

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H102").Value = 2676.6738
$ws.Range("I102").Value = 2518.8684
$ws.Range("J102").Value = 3426.25
$ws.Range("K102").Value = 2518.8684
$ws.Range("L102").Value = 3426.25
$ws.Range("M102").Value = -896.8683999999998
$ws.Range("N102").Value = -6670.25

$ws.Range("H113").Value = 793.8333
$ws.Range("I113").Value = 762.6
$ws.Range("K113").Value = 762.6
$ws.Range("M113").Value = 1407.4

$ws.Range("H126").Value = 5892.1055
$ws.Range("I126").Value = 3659.6
$ws.Range("J126").Value = 8372.666999999999
$ws.Range("K126").Value = 10978.8
$ws.Range("L126").Value = 25118.001
$ws.Range("M126").Value = -8508.799999999999
$ws.Range("N126").Value = -30058.001

$ws.Range("H132").Value = 4935.643
$ws.Range("I132").Value = 5091.5835
$ws.Range("K132").Value = 15274.7505
$ws.Range("M132").Value = -12744.7505

$ws.Range("H140").Value = 78926.664
$ws.Range("J140").Value = 78926.664
$ws.Range("L140").Value = 78926.664
$ws.Range("N140").Value = -89286.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12509995
$ws.Range("J20").Value = 12509995
$ws.Range("L20").Value = 12509995
$ws.Range("N20").Value = -12510447

$ws.Range("H40").Value = 2082.8
$ws.Range("I40").Value = 1844.4
$ws.Range("J40").Value = 2798
$ws.Range("K40").Value = 1844.4
$ws.Range("L40").Value = 2798
$ws.Range("M40").Value = -1708.4
$ws.Range("N40").Value = -3070

$ws.Range("H46").Value = 2308.6365
$ws.Range("J46").Value = 2737.125
$ws.Range("L46").Value = 2737.125
$ws.Range("N46").Value = -3113.125

$ws.Range("H61").Value = 2969.8333
$ws.Range("I61").Value = 3059.111
$ws.Range("J61").Value = 2702
$ws.Range("K61").Value = 3059.111
$ws.Range("L61").Value = 2702
$ws.Range("M61").Value = -2857.111
$ws.Range("N61").Value = -3106

$ws.Range("H68").Value = 2520.3635
$ws.Range("I68").Value = 2512.6667
$ws.Range("K68").Value = 2512.6667
$ws.Range("M68").Value = -1763.6667

$ws.Range("H71").Value = 2520.3635
$ws.Range("I71").Value = 2512.6667
$ws.Range("K71").Value = 12563.3335
$ws.Range("M71").Value = -8819.333500000001

$ws.Range("H82").Value = 4643.8823
$ws.Range("I82").Value = 1743.2222
$ws.Range("J82").Value = 7907.125
$ws.Range("K82").Value = 1743.2222
$ws.Range("L82").Value = 7907.125
$ws.Range("M82").Value = -1382.2222
$ws.Range("N82").Value = -8629.125

$ws.Range("H85").Value = 4643.8823
$ws.Range("I85").Value = 1743.2222
$ws.Range("J85").Value = 7907.125
$ws.Range("K85").Value = 1743.2222
$ws.Range("L85").Value = 7907.125
$ws.Range("M85").Value = -495.2221999999999
$ws.Range("N85").Value = -10403.125

$ws.Range("H93").Value = 697.3
$ws.Range("I93").Value = 362.33334
$ws.Range("K93").Value = 362.33334
$ws.Range("M93").Value = 885.66666

$ws.Range("H100").Value = 1783.5
$ws.Range("I100").Value = 1724.25
$ws.Range("J100").Value = 1902
$ws.Range("K100").Value = 1724.25
$ws.Range("L100").Value = 1902
$ws.Range("M100").Value = -1183.25
$ws.Range("N100").Value = -2984

$ws.Range("H113").Value = 2969.8333
$ws.Range("I113").Value = 3059.111
$ws.Range("J113").Value = 2702
$ws.Range("K113").Value = 3059.111
$ws.Range("L113").Value = 2702
$ws.Range("M113").Value = -889.1109999999999
$ws.Range("N113").Value = -7042

$ws.Range("H122").Value = 3997.7778
$ws.Range("I122").Value = 3998
$ws.Range("J122").Value = 3997.3333
$ws.Range("K122").Value = 11994
$ws.Range("L122").Value = 11991.9999
$ws.Range("M122").Value = -9544
$ws.Range("N122").Value = -16891.9999

$ws.Range("H133").Value = 84332.664
$ws.Range("J133").Value = 84332.664
$ws.Range("L133").Value = 84332.664
$ws.Range("N133").Value = -89392.664

$ws.Range("H134").Value = 25000
$ws.Range("J134").Value = 25000
$ws.Range("L134").Value = 25000
$ws.Range("N134").Value = -35140

$ws.Range("H136").Value = 3126.739
$ws.Range("I136").Value = 2181
$ws.Range("K136").Value = 6543
$ws.Range("M136").Value = -3993

$ws.Range("H139").Value = 66851.81
$ws.Range("J139").Value = 66851.81
$ws.Range("L139").Value = 66851.81
$ws.Range("N139").Value = -77131.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H41").Value = 72798.60000000001
$ws.Range("I41").Value = 69998
$ws.Range("J41").Value = 73498.75
$ws.Range("K41").Value = 69998
$ws.Range("L41").Value = 73498.75
$ws.Range("M41").Value = -69608
$ws.Range("N41").Value = -74278.75

$ws.Range("H62").Value = 9999
$ws.Range("J62").Value = 9999
$ws.Range("L62").Value = 9999
$ws.Range("N62").Value = -11247

$ws.Range("H65").Value = 9999
$ws.Range("J65").Value = 9999
$ws.Range("L65").Value = 49995
$ws.Range("N65").Value = -56235

$ws.Range("H100").Value = 2145.1538
$ws.Range("I100").Value = 386.75
$ws.Range("K100").Value = 773.5
$ws.Range("M100").Value = -232.5

$ws.Range("H107").Value = 5230.273
$ws.Range("I107").Value = 7258.067
$ws.Range("J107").Value = 885
$ws.Range("K107").Value = 21774.201
$ws.Range("L107").Value = 2655
$ws.Range("M107").Value = -19854.201
$ws.Range("N107").Value = -6495

$ws.Range("H113").Value = 910150.9399999999
$ws.Range("I113").Value = 1112078.9
$ws.Range("K113").Value = 3336236.7
$ws.Range("M113").Value = -3334066.7

$ws.Range("H132").Value = 3893.7585
$ws.Range("I132").Value = 3246.6667
$ws.Range("K132").Value = 9740.000100000001
$ws.Range("M132").Value = -7210.000100000001

$ws.Range("H136").Value = 14499.4
$ws.Range("I136").Value = 13998.538
$ws.Range("K136").Value = 41995.614
$ws.Range("M136").Value = -39445.614
